$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptos list values (Price and Volume(1h) columns).
# Cells whose new text is a plain decimal number (e.g. "42.68") must be
# forced to Text so Excel keeps the literal formatted string instead of
# silently converting the cell to a numeric value.

$ws.Range('D2').Value = '30.256.09'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '1.991.38'
$ws.Range('E3').Value = '  +6.09%  '
$ws.Range('E4').Value = '  -0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '325.37'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.78%  '
$ws.Range('E6').Value = '  -0.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5096'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.14%  '
$ws.Range('E8').Value = '  +4.91%  '
$ws.Range('E9').Value = '  +6.00%  '
$ws.Range('E10').Value = '  +3.54%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.68'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.20'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +2.56%  '
$ws.Range('D13').Value = '1.994.44'
$ws.Range('E13').Value = '  +6.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.495'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +3.21%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.403'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9993'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '94.03'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.50%  '
$ws.Range('E18').Value = '  +2.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06555'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.44%  '
$ws.Range('E20').Value = '  +4.70%  '
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.074'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +4.05%  '
$ws.Range('D23').Value = '30.321.88'
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.59'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +4.01%  '
$ws.Range('E25').Value = '  +2.28%  '
$ws.Range('D26').Value = '2.222.90'
$ws.Range('E26').Value = '  +6.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.55'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +6.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.29'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.381'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +6.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '130.80'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.132'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +5.50%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1052'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.66%  '
$ws.Range('E33').Value = '  +2.33%  '
$ws.Range('E34').Value = '  +3.03%  '
$ws.Range('E35').Value = '  +12.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02478'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.17%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.391'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +2.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06523'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2194'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.03%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.910'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +4.86%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6576'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.65%  '
$ws.Range('E42').Value = '  +4.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.225'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.65'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6126'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.196'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +5.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.669'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.21'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.226'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '79.40'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06874'
$ws.Range('D51').ClearFormats()
